$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M3").Value = 1.06
$ws.Range("N3").Value = 8
$ws.Range("Q3").Value = 2.08
$ws.Range("R3").Value = 1.73
$ws.Range("G5").Value = 2.3
$ws.Range("H5").Value = 3.3
$ws.Range("J5").Value = 2.88
$ws.Range("K5").Value = 2.2
$ws.Range("L5").Value = 3.25
$ws.Range("O5").Value = 1.25
$ws.Range("P5").Value = 3.75
$ws.Range("Q5").Value = 1.85
$ws.Range("R5").Value = 1.95
$ws.Range("S5").Value = 1.36
$ws.Range("T5").Value = 3
$ws.Range("U5").Value = 1.67
$ws.Range("V5").Value = 2.1
$ws.Range("W5").Value = 9
$ws.Range("Y5").Value = 9.5
$ws.Range("Z5").Value = 21
$ws.Range("AB5").Value = 26
$ws.Range("AC5").Value = 11
$ws.Range("AM5").Value = 29
$ws.Range("AP5").Value = 21
$ws.Range("AR5").Value = 51
$ws.Range("AT5").Value = 3
$ws.Range("AU5").Value = 7.5
$ws.Range("AY5").Value = 23
$ws.Range("G6").Value = 2.75
$ws.Range("H6").Value = 3.2
$ws.Range("I6").Value = 2.3
$ws.Range("J6").Value = 3.6
$ws.Range("L6").Value = 3.2
$ws.Range("M6").Value = 1.07
$ws.Range("N6").Value = 8.5
$ws.Range("W6").Value = 8
$ws.Range("X6").Value = 13
$ws.Range("Z6").Value = 29
$ws.Range("AC6").Value = 8.5
$ws.Range("AH6").Value = 7.5
$ws.Range("AI6").Value = 11
$ws.Range("AK6").Value = 23
$ws.Range("AL6").Value = 21
$ws.Range("AM6").Value = 34
$ws.Range("AN6").Value = 4.75
$ws.Range("AO6").Value = 17
$ws.Range("AQ6").Value = 51
$ws.Range("AX6").Value = 15
$ws.Range("AY6").Value = 26
$ws.Range("AZ6").Value = 51
